$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the current B:D values (rows 1-5) before touching anything.
$data = @()
for ($r = 1; $r -le 5; $r++) {
    $row = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2
    )
    $data += ,$row
}

# Clear the old B:D block now that the values are captured.
$ws.Range("B1:D5").Clear()

# Write the values shifted one column to the left: B->A, C->B, D->C.
for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Fix the header typo: "Edad" -> "Edada1"
$ws.Cells.Item(1, 3).Value = "Edada1"

# Add the block of styled (underlined) empty cells I14:L18
$fmtRange = $ws.Range("I14:L18")
$fmtRange.Font.Underline = 1

# Update the active selection to F3
$ws.Range("F3").Select() | Out-Null
